# Updates the cryptos list (Coin / Link / Price / Volume(1h) columns) to
# match the "Updated cryptos list ... with GitHub Actions" commit.
#
# Price/Volume values are numeric-looking strings (e.g. "531.50", "1.00")
# that Excel would otherwise silently coerce to numbers (losing the
# trailing zero / exact text) when assigned through .Value. Force them to
# stay text by flipping the cell to the "@" (Text) number format for the
# assignment, then restore the "Normal" style so no stray style index is
# left behind on the cell.
function Set-TextValue {
    param($Cell, $Val)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Val
    $Cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose Coin/Link stay put - only Price (D) and/or Volume (E) change.
$rows = @(
    @{ Row = 2;  D = "58.730.19";  E = "  -5.10%  " },
    @{ Row = 3;  D = "2.484.83";   E = "  -3.53%  " },
    @{ Row = 4;  D = $null;        E = "  -0.08%  " },
    @{ Row = 5;  D = "531.50";     E = "  -3.21%  " },
    @{ Row = 6;  D = "142.53";     E = $null },
    @{ Row = 7;  D = "0.996";      E = "  -0.30%  " },
    @{ Row = 8;  D = $null;        E = "  -4.52%  " },
    @{ Row = 9;  D = "2.510.91";   E = "  -2.77%  " },
    @{ Row = 10; D = $null;        E = "  -4.63%  " },
    @{ Row = 11; D = $null;        E = "  -2.81%  " },
    @{ Row = 12; D = "5.48";       E = "  +0.83%  " },
    @{ Row = 13; D = "0.348";      E = "  -4.52%  " },
    @{ Row = 14; D = "2.920.79";   E = "  -3.69%  " },
    @{ Row = 15; D = "23.60";      E = "  -7.74%  " },
    @{ Row = 16; D = "58.595.20";  E = "  -5.16%  " },
    @{ Row = 17; D = $null;        E = "  -4.81%  " },
    @{ Row = 18; D = "2.496.56";   E = "  -3.26%  " },
    @{ Row = 19; D = "11.28";      E = "  -2.14%  " },
    @{ Row = 20; D = "4.24";       E = "  -6.57%  " },
    @{ Row = 21; D = "320.32";     E = "  -5.08%  " },
    @{ Row = 22; D = $null;        E = "  -0.06%  " },
    @{ Row = 23; D = "5.71";       E = "  -5.50%  " },
    @{ Row = 24; D = "60.66";      E = "  -4.30%  " },
    @{ Row = 25; D = $null;        E = "  -11.39%  " },
    @{ Row = 34; D = $null;        E = "  -0.27%  " },
    @{ Row = 35; D = "156.63";     E = "  -2.62%  " },
    @{ Row = 36; D = "1.39";       E = "  -1.61%  " },
    @{ Row = 37; D = "18.39";      E = "  -4.22%  " },
    @{ Row = 38; D = "4.33";       E = "  -9.91%  " },
    @{ Row = 39; D = "1.60";       E = "  -10.96%  " },
    @{ Row = 40; D = "5.90";       E = "  -1.10%  " },
    @{ Row = 41; D = "307.29";     E = "  -7.26%  " },
    @{ Row = 42; D = "36.67";      E = "  -2.31%  " },
    @{ Row = 43; D = $null;        E = "  -7.82%  " },
    @{ Row = 44; D = "0.784";      E = "  -14.65%  " },
    @{ Row = 45; D = "0.994";      E = "  -0.34%  " },
    @{ Row = 46; D = "0.592";      E = "  -2.27%  " },
    @{ Row = 47; D = $null;        E = "  -1.57%  " },
    @{ Row = 48; D = "124.41";     E = "  +0.88%  " },
    @{ Row = 49; D = "0.0919";     E = "  -4.71%  " },
    @{ Row = 50; D = "18.44";      E = "  -5.71%  " },
    @{ Row = 51; D = "0.0511";     E = "  -6.52%  " }
)

foreach ($r in $rows) {
    if ($null -ne $r.D) {
        Set-TextValue $ws.Cells.Item($r.Row, 4) $r.D
    }
    if ($null -ne $r.E) {
        Set-TextValue $ws.Cells.Item($r.Row, 5) $r.E
    }
}

# Rows 26-33 were fully re-sorted (Coin/Link/Price/Volume all change in place).
$reorder = @(
    @{ Row = 26; B = "Binance-PegBSC-USD";        C = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd";   D = "0.995";    E = "  -0.49%  " },
    @{ Row = 27; B = "Kaspa";                      C = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas";                   D = "0.161";    E = "  -3.98%  " },
    @{ Row = 28; B = "WrappedeETH";                C = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth";           D = "2.597.90"; E = "  -3.65%  " },
    @{ Row = 29; B = "InternetComputer(DFINITY)";  C = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; D = "7.67";     E = "  -5.62%  " },
    @{ Row = 30; B = "Aptos";                      C = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt";                   D = "6.77";     E = "  -9.14%  " },
    @{ Row = 31; B = "Fetch.AI";                   C = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet";             D = "1.24";     E = "  -6.84%  " },
    @{ Row = 32; B = "PEPE";                       C = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe";                   D = "0.0₃0758"; E = "  -9.58%  " },
    @{ Row = 33; B = "PancakeSwap";                C = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake";            D = "1.77";     E = "  -6.74%  " }
)

foreach ($r in $reorder) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    Set-TextValue $ws.Cells.Item($r.Row, 4) $r.D
    Set-TextValue $ws.Cells.Item($r.Row, 5) $r.E
}
